$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells keep their original text representation (the source workbook
# stores these as inline strings, e.g. "528.91" or "0.999", not numbers), so
# force text number-format before assigning the value to avoid Excel silently
# converting numeric-looking strings into real numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.039.33"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.18%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.230.85"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.09%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.28%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "528.91"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +3.96%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "171.04"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.68%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +1.80%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.17%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.227.69"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.18%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.605"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.18%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "52.96"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -5.18%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.133"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +4.49%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000254"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +2.03%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.90%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.745.27"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.28%  "
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.58%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.231.38"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.22%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "62.906.76"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.16%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.17"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.63%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.04"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +3.88%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +4.20%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "365.80"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.73%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.74"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +4.53%  "
$ws.Range("B24").NumberFormat = "@"
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").NumberFormat = "@"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "80.78"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +3.11%  "
$ws.Range("B25").NumberFormat = "@"
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").NumberFormat = "@"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.03"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +9.06%  "
$ws.Range("B26").NumberFormat = "@"
$ws.Range("B26").Value = "RenderToken"
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.11"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +3.71%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +2.22%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.96%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "11.20"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.58%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.19"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.13%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "28.40"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +2.25%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "633.67"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.07%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.40"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.25%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.19"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +2.36%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.106"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +4.49%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "56.88"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -2.90%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.07%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "36.41"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +3.70%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.39%  "
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = "PEPE"
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0₃0724"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +14.07%  "
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "FirstDigitalUSD"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.00%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +2.93%  "
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.868.01"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.64%  "
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = "Fetch.AI"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.54"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +10.23%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.93"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +6.99%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +4.68%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +7.74%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +5.44%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.99%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +2.75%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "133.71"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +2.65%  "
